$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.929.65"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "3.786.65"
$ws.Range("E3").Value = "  +3.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "428.53"
$ws.Range("E5").Value = "  +6.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.99"
$ws.Range("E6").Value = "  +10.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  +4.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +4.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -4.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000312"
$ws.Range("E11").Value = "  -7.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.68"
$ws.Range("E12").Value = "  +7.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.48"
$ws.Range("E13").Value = "  +11.60%  "

$ws.Range("D14").Value = "4.381.66"
$ws.Range("E14").Value = "  +4.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("E15").Value = "  +3.99%  "

$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "3.782.10"
$ws.Range("E17").Value = "  +3.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.94"
$ws.Range("E18").Value = "  +4.64%  "

$ws.Range("E19").Value = "  +7.95%  "

$ws.Range("D20").Value = "66.076.55"
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "407.56"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.23"
$ws.Range("E22").Value = "  +7.16%  "

$ws.Range("E23").Value = "  +7.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.91"
$ws.Range("E24").Value = "  +2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.73"
$ws.Range("E25").Value = "  +4.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.28"
$ws.Range("E26").Value = "  +8.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.71"
$ws.Range("E27").Value = "  +39.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  +10.70%  "

$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.81"
$ws.Range("E30").Value = "  +12.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "711.55"
$ws.Range("E31").Value = "  +4.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("E32").Value = "  +16.47%  "

$ws.Range("E33").Value = "  +4.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.03"
$ws.Range("E34").Value = "  +8.92%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("E37").Value = "  +35.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.60"
$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("E39").Value = "  +6.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  +45.57%  "

$ws.Range("E41").Value = "  +6.75%  "

$ws.Range("D42").Value = "0.0₃0686"
$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("E43").Value = "  +6.96%  "

$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.325"
$ws.Range("E45").Value = "  +14.46%  "

$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +7.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.15"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("E48").Value = "  +3.43%  "

$ws.Range("E49").Value = "  +4.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.58"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  +4.45%  "
